$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.276.93"
$ws.Range("E2").Value = '  +1.91%  '

$ws.Range("D3").Value = "'1.913.64"
$ws.Range("E3").Value = '  +2.20%  '

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = "'328.29"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").Value = "'0.4632"
$ws.Range("E7").Value = '  +1.05%  '

$ws.Range("D8").Value = "'0.3946"
$ws.Range("E8").Value = '  +2.50%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = "'0.07950"
$ws.Range("E9").Value = '  +1.19%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = '  +1.00%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = "'22.33"
$ws.Range("E11").Value = '  +2.87%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'1.937.25"
$ws.Range("E12").Value = '  +2.80%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = "'7.107"
$ws.Range("E13").Value = '  +2.04%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'5.774"
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = "'0.06963"
$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = "'88.62"
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.00001009"
$ws.Range("E18").Value = '  +0.52%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = "'17.16"
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = "'29.276.98"
$ws.Range("E21").Value = '  +1.87%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = "'5.360"
$ws.Range("E22").Value = '  +1.67%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = "'11.11"
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = "'2.153.84"
$ws.Range("E24").Value = '  +2.03%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'2.060"
$ws.Range("E25").Value = '  -3.06%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'157.09"
$ws.Range("E26").Value = '  +2.49%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'19.51"
$ws.Range("E27").Value = '  +1.46%  '

$ws.Range("D28").Value = "'6.165"
$ws.Range("E28").Value = '  +6.96%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = "'1.994"
$ws.Range("E29").Value = '  +2.53%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = "'118.84"
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.09376"
$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'0.9247"
$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'5.360"
$ws.Range("E33").Value = '  +1.26%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = "'1.357"
$ws.Range("E34").Value = '  +1.77%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'3.281"
$ws.Range("E35").Value = '  -0.86%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = "'1.195"
$ws.Range("E36").Value = '  +4.14%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.05839"
$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.02113"
$ws.Range("E38").Value = '  +1.92%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'7.970"
$ws.Range("E39").Value = '  +3.61%  '

$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = "'1.003"
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.5755"
$ws.Range("E41").Value = '  +2.44%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = "'0.1802"
$ws.Range("E42").Value = '  +0.93%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'9.948"
$ws.Range("E43").Value = '  +1.02%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'2.298"
$ws.Range("E44").Value = '  +8.19%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'12.04"
$ws.Range("E45").Value = '  +2.26%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.5429"
$ws.Range("E46").Value = '  +3.07%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = "'0.07072"
$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'1.881"
$ws.Range("E48").Value = '  +3.10%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = "'2.564"
$ws.Range("E49").Value = '  +6.19%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'113.57"
$ws.Range("E50").Value = '  +0.27%  '

$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = "'1.073"
$ws.Range("E51").Value = '  -3.83%  '

$ws.Range("D2:D51").Style = "Normal"
